# B6-PowerPoint.pptx edit script
#
# 1) Three tables (on slides 14, 15 and 16) get their table style switched
#    from the default "Table_0" style ({CFA07F61-4F73-4064-B370-DC8941109A24})
#    to the built-in "Medium Style 2 - Accent 1" style
#    ({C8ECCB3E-1F9C-4C88-B99B-70F23EFB9265}).
#
# 2) The deck's theme (design) is switched from the pink/red-violet
#    "Integral" theme to the standard blue/grey "Office Theme" -- i.e. every
#    theme colour slot on the slide master's theme is repointed at the
#    stock Office colour scheme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the three tables.
# ---------------------------------------------------------------------
$newTableStyle = "{C8ECCB3E-1F9C-4C88-B99B-70F23EFB9265}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the "Integral" / Red Violet theme colours for the stock
#    "Office Theme" / Office colours on the presentation's theme.
# ---------------------------------------------------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# RGB() is expressed as a plain 0x00BBGGRR COM colour integer since the
# RGB() helper isn't available in this host.
$tcs.Item(1).RGB  = 0x000000    # Dark 1    -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF    # Light 1   -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444    # Dark 2    -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7    # Light 2   -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B    # Accent 1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED    # Accent 2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5    # Accent 3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF    # Accent 4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244    # Accent 5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70    # Accent 6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305    # Hyperlink -> 0563C1
$tcs.Item(12).RGB = 0x724F95    # Followed Hyperlink -> 954F72
